$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.861.82'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.03%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.808.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.81%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '353.93'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.20%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '112.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.36%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.556'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.32%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.625'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.32%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.40'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E11').Value = '  -0.10%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0838'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.63%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.96'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.78%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.78'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.28%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.246.99'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.49%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.792.87'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.46%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.945'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.94%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.782.15'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.93%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.84%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +7.40%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.33%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0973'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.75%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.33'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.25%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '267.72'
$ws.Range('D24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.52%  '

$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.15%  '

$ws.Range('E28').Value = '  +0.23%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '39.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +14.62%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '10.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.54%  '

$ws.Range('B31').Value = 'OKB'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '52.58'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.48%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.15'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.29%  '

$ws.Range('B33').Value = 'Toncoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.12'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.24%  '

$ws.Range('E34').Value = '  +3.75%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0895'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.11%  '

$ws.Range('E36').Value = '  +8.78%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.998'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.21%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.93'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.20%  '

$ws.Range('E39').Value = '  +2.33%  '

$ws.Range('E40').Value = '  +5.00%  '

$ws.Range('E41').Value = '  +2.86%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.53'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.90%  '

$ws.Range('E43').Value = '  +1.51%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '120.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.95%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.11%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.48'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +10.21%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.47'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.73%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.108.69'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.52%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.965'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.04%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.76%  '

$ws.Range('E51').Value = '  +8.47%  '
